$d = $word.ActiveDocument

# Update the date heading (first paragraph).
# Find.Execute locates the match in place (Replace = 0, wdReplaceNone) and
# then we assign the matched Range's .Text directly so the edit stays
# confined to that exact Range. (Using Find's built-in auto-replace here
# would rewrite every occurrence of the search text anywhere in the
# document, which is wrong when the same expression text recurs in
# multiple table cells.)
$titleRange = $d.Paragraphs.Item(1).Range
$null = $titleRange.Find.Execute("2024-12-27 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$titleRange.Text = "2024-12-28 Saturday"

# Update each answer cell in the table (scoped per-cell via the located Range).
$t = $d.Tables.Item(1)
$c = $t.Cell(1,1).Range
$null = $c.Find.Execute("69÷2=34, 1", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$c.Text = "53÷7=7, 4"
$c = $t.Cell(1,2).Range
$null = $c.Find.Execute("22÷7=3, 1", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$c.Text = "81÷9=9, 0"
$c = $t.Cell(1,3).Range
$null = $c.Find.Execute("51÷7=7, 2", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$c.Text = "99÷5=19, 4"
$c = $t.Cell(1,4).Range
$null = $c.Find.Execute("50÷7=7, 1", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$c.Text = "54÷2=27, 0"
$c = $t.Cell(1,5).Range
$null = $c.Find.Execute("85÷9=9, 4", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$c.Text = "52÷6=8, 4"
$c = $t.Cell(5,1).Range
$null = $c.Find.Execute("58÷5=11, 3", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$c.Text = "80÷5=16, 0"
$c = $t.Cell(5,2).Range
$null = $c.Find.Execute("52÷6=8, 4", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$c.Text = "60÷5=12, 0"
$c = $t.Cell(5,3).Range
$null = $c.Find.Execute("46÷6=7, 4", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$c.Text = "93÷5=18, 3"
$c = $t.Cell(5,4).Range
$null = $c.Find.Execute("78÷2=39, 0", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$c.Text = "48÷2=24, 0"
$c = $t.Cell(5,5).Range
$null = $c.Find.Execute("99÷4=24, 3", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$c.Text = "77÷9=8, 5"
$c = $t.Cell(9,1).Range
$null = $c.Find.Execute("34÷6=5, 4", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$c.Text = "52÷6=8, 4"
$c = $t.Cell(9,2).Range
$null = $c.Find.Execute("27÷4=6, 3", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$c.Text = "49÷8=6, 1"
$c = $t.Cell(9,3).Range
$null = $c.Find.Execute("76÷4=19, 0", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$c.Text = "15÷6=2, 3"
$c = $t.Cell(9,4).Range
$null = $c.Find.Execute("99÷9=11, 0", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$c.Text = "70÷3=23, 1"
$c = $t.Cell(9,5).Range
$null = $c.Find.Execute("41÷8=5, 1", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$c.Text = "37÷6=6, 1"
$c = $t.Cell(13,1).Range
$null = $c.Find.Execute("48÷3=16, 0", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$c.Text = "84÷9=9, 3"
$c = $t.Cell(13,2).Range
$null = $c.Find.Execute("28÷8=3, 4", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$c.Text = "25÷4=6, 1"
$c = $t.Cell(13,3).Range
$null = $c.Find.Execute("84÷3=28, 0", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$c.Text = "53÷7=7, 4"
$c = $t.Cell(13,4).Range
$null = $c.Find.Execute("24÷8=3, 0", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$c.Text = "83÷4=20, 3"
$c = $t.Cell(13,5).Range
$null = $c.Find.Execute("76÷2=38, 0", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$c.Text = "76÷4=19, 0"
$c = $t.Cell(17,1).Range
$null = $c.Find.Execute("66÷8=8, 2", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$c.Text = "56÷7=8, 0"
$c = $t.Cell(17,2).Range
$null = $c.Find.Execute("34÷2=17, 0", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$c.Text = "60÷3=20, 0"
$c = $t.Cell(17,3).Range
$null = $c.Find.Execute("81÷6=13, 3", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$c.Text = "32÷5=6, 2"
$c = $t.Cell(17,4).Range
$null = $c.Find.Execute("73÷9=8, 1", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$c.Text = "88÷9=9, 7"
$c = $t.Cell(17,5).Range
$null = $c.Find.Execute("75÷4=18, 3", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$c.Text = "25÷6=4, 1"
